$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The "R30" rule row (row 10) had its "From" value (column C) changed
# from 18 to 1.
$ws.Range("C10").Value = 1
